$wb = $excel.ActiveWorkbook

$newPeriod = "31/01/2022 - 06/02/2022"

# ---------------------------------------------------------------------------
# Sheet "Asl Sorveglianza" (sheet1): add week block 31/01/2022 - 06/02/2022
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Asl Sorveglianza")

$ws1.Range("A112").Value = $newPeriod
$ws1.Range("B112").Value = "AZIENDA USL TOSCANA SUD-EST"
$ws1.Range("C112").Value = 38

$ws1.Range("A113").Value = $newPeriod
$ws1.Range("B113").Value = "AZIENDA USL TOSCANA CENTRO"
$ws1.Range("C113").Value = 40

$ws1.Range("A114").Value = $newPeriod
$ws1.Range("B114").Value = "AZIENDA USL TOSCANA NORD-OVEST"
$ws1.Range("B114").Font.Color = 0
$ws1.Range("C114").Value = 90

$ws1.Range("B115").Value = "Totale"
$ws1.Range("C115").Value = 169

$ws1.Range("A112").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Professione" (sheet2): add week block 31/01/2022 - 06/02/2022
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Professione")

$ws2.Range("A82").Value = $newPeriod
$ws2.Range("B82").Value = "Insegnante"
$ws2.Range("C82").Value = 161
$ws2.Range("D82").Value = 1280
$ws2.Range("D82").NumberFormat = "#,##0"

$ws2.Range("A83").Value = $newPeriod
$ws2.Range("B83").Value = "Personale non docente"
$ws2.Range("C83").Value = 8
$ws2.Range("D83").Value = 30
$ws2.Range("D83").NumberFormat = "#,##0"

$ws2.Range("B84").Value = "Totale"
$ws2.Range("C84").Value = 169
$ws2.Range("D84").Value = 1310
$ws2.Range("D84").NumberFormat = "#,##0"

$ws2.Range("A82").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Sesso ed età" (sheet3): add week block 31/01/2022 - 06/02/2022
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sesso ed età")

$ws3.Range("A150").Value = $newPeriod
$ws3.Range("B150").Value = "0-18"
$ws3.Range("C150").Value = "F"
$ws3.Range("D150").Value = 1

$ws3.Range("A151").Value = $newPeriod
$ws3.Range("B151").Value = "19-34"
$ws3.Range("C151").Value = "F"
$ws3.Range("D151").Value = 18

$ws3.Range("A152").Value = $newPeriod
$ws3.Range("B152").Value = "19-34"
$ws3.Range("B152").Font.Color = 0
$ws3.Range("C152").Value = "M"
$ws3.Range("D152").Value = 3

$ws3.Range("A153").Value = $newPeriod
$ws3.Range("B153").Value = "35-49"
$ws3.Range("C153").Value = "F"
$ws3.Range("D153").Value = 75

$ws3.Range("A154").Value = $newPeriod
$ws3.Range("B154").Value = "35-49"
$ws3.Range("C154").Value = "M"
$ws3.Range("D154").Value = 10

$ws3.Range("A155").Value = $newPeriod
$ws3.Range("B155").Value = "50-64"
$ws3.Range("C155").Value = "F"
$ws3.Range("D155").Value = 52

$ws3.Range("A156").Value = $newPeriod
$ws3.Range("B156").Value = "50-64"
$ws3.Range("C156").Value = "M"
$ws3.Range("D156").Value = 7

$ws3.Range("A157").Value = $newPeriod
$ws3.Range("B157").Value = "65-79"
$ws3.Range("C157").Value = "F"
$ws3.Range("D157").Value = 3

$ws3.Range("A158:D158").Select() | Out-Null

$ws3.Activate() | Out-Null
